# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# A fonte de dados foi re-sincronizada: alguns jogos trocaram de posição
# (mesma "id" de linha na coluna A, mas os dados do jogo vieram trocados
# entre duas linhas consecutivas) e algumas odds de jogos futuros foram
# atualizadas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB, $firstCol, $lastCol) {
    $valsA = @{}
    $valsB = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $valsA[$c] = $sheet.Cells.Item($rowA, $c).Value2
        $valsB[$c] = $sheet.Cells.Item($rowB, $c).Value2
    }
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $sheet.Cells.Item($rowA, $c).Value = $valsB[$c]
        $sheet.Cells.Item($rowB, $c).Value = $valsA[$c]
    }
}

# Column B (id) through AC (PL_AhUnder) swap, column A (row index) is left untouched.
Swap-Rows $ws 142 143 2 29
Swap-Rows $ws 167 168 2 29
Swap-Rows $ws 222 223 2 29
Swap-Rows $ws 241 242 2 29

# Odds updates for upcoming fixtures (not yet played).
$ws.Range("O245").Value = 3.3
$ws.Range("P245").Value = 3
$ws.Range("U245").Value = 2.05
$ws.Range("V245").Value = 1.8

$ws.Range("N252").Value = 2.8
$ws.Range("P252").Value = 2.2
$ws.Range("R252").Value = 1.825
$ws.Range("S252").Value = 2.025

$ws.Range("N253").Value = 2.375
$ws.Range("P253").Value = 2.625
$ws.Range("Q253").Value = 0
$ws.Range("R253").Value = 1.8
$ws.Range("S253").Value = 2.05
$ws.Range("U253").Value = 1.95
$ws.Range("V253").Value = 1.9
